$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# A new entry was logged on 2024-09-20 15:37:11 ("check the loan yo").
# It belongs at the top of the "Others" group's September data, so every
# existing data row from 45 down to the end of the sheet shifts down by
# one row; insert a fresh row 45 and fill in the new entry there.
$ws.Rows.Item(45).Insert()
$ws.Range("R45").Value = "check the loan yo"
$ws.Range("S45").Value = "2024-09-20 15:37:11"
